$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5664
$ws.Range("J3").Value = 6040
$ws.Range("C4").Value = 1836
$ws.Range("J4").Value = 1308
$ws.Range("J5").Value = 463
$ws.Range("J6").Value = 7723
$ws.Range("C7").Value = 28380
$ws.Range("J7").Value = 21198

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 193
$ws.Range("J7").Value = 297

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 365
$ws.Range("J3").Value = 405
$ws.Range("J4").Value = 77
$ws.Range("J6").Value = 449
$ws.Range("J7").Value = 1331

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 127
$ws.Range("J7").Value = 428

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 323
$ws.Range("J6").Value = 339
$ws.Range("J7").Value = 981

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 196
$ws.Range("J7").Value = 658

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J6").Value = 189
$ws.Range("J7").Value = 538

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 94
$ws.Range("J7").Value = 333

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 169
$ws.Range("J4").Value = 93
$ws.Range("J7").Value = 621
$ws.Range("J8").Value = 1331
$ws.Range("J10").Value = 145
$ws.Range("J15").Value = 233
$ws.Range("J16").Value = 82
$ws.Range("J19").Value = 617
$ws.Range("J20").Value = 439
$ws.Range("J22").Value = 55
$ws.Range("J24").Value = 65
$ws.Range("J26").Value = 46
$ws.Range("J27").Value = 128
$ws.Range("J29").Value = 1185
$ws.Range("J31").Value = 192
$ws.Range("J32").Value = 35
$ws.Range("J33").Value = 981
$ws.Range("J34").Value = 101
$ws.Range("J37").Value = 658
$ws.Range("J41").Value = 133
$ws.Range("J42").Value = 887
$ws.Range("J48").Value = 248
$ws.Range("J49").Value = 144
$ws.Range("C52").Value = 647
$ws.Range("J53").Value = 297
$ws.Range("J54").Value = 410
$ws.Range("J55").Value = 284
$ws.Range("J57").Value = 92
$ws.Range("J63").Value = 75
$ws.Range("J64").Value = 140
$ws.Range("J65").Value = 538
$ws.Range("J67").Value = 805
$ws.Range("J76").Value = 311
$ws.Range("J78").Value = 263
$ws.Range("J79").Value = 605
$ws.Range("J80").Value = 33
$ws.Range("J83").Value = 428
$ws.Range("J85").Value = 884
$ws.Range("J87").Value = 73
$ws.Range("J88").Value = 225
$ws.Range("J89").Value = 281
$ws.Range("J91").Value = 238
$ws.Range("J92").Value = 64
$ws.Range("J94").Value = 213
$ws.Range("J96").Value = 248
$ws.Range("J97").Value = 175
$ws.Range("J98").Value = 155
$ws.Range("J99").Value = 333
$ws.Range("C101").Value = 28380
$ws.Range("J101").Value = 21198

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J7").Value = 621
$ws.Range("J3").Value = 187
$ws.Range("J6").Value = 200

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 74
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 192

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 202
$ws.Range("J3").Value = 306
$ws.Range("J6").Value = 214
$ws.Range("J7").Value = 805

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 196
$ws.Range("J7").Value = 410

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 411
$ws.Range("J6").Value = 308
$ws.Range("J7").Value = 1185

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 39
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 248

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 153
$ws.Range("J3").Value = 180
$ws.Range("J6").Value = 231
$ws.Range("J7").Value = 617

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J3").Value = 65
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J6").Value = 72
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 193
$ws.Range("J3").Value = 178
$ws.Range("J6").Value = 459
$ws.Range("J7").Value = 887

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 28
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 263

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 63
$ws.Range("J7").Value = 284

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J2").Value = 20
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 74
$ws.Range("J3").Value = 67
$ws.Range("J7").Value = 248

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 69
$ws.Range("J3").Value = 96
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 238

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 170
$ws.Range("J7").Value = 605

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 140

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 121
$ws.Range("J5").Value = 10
$ws.Range("J6").Value = 116
$ws.Range("J7").Value = 439

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 30
$ws.Range("J7").Value = 101

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 117
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J5").Value = 4
$ws.Range("J6").Value = 97
$ws.Range("J7").Value = 233

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J2").Value = 28
$ws.Range("J6").Value = 96
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 48
$ws.Range("J7").Value = 169

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J2").Value = 15
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 225

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 78
$ws.Range("J7").Value = 281

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J5").Value = 18
$ws.Range("J6").Value = 259
$ws.Range("J7").Value = 884

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 163
$ws.Range("C4").Value = 40
$ws.Range("J4").Value = 19
$ws.Range("C7").Value = 647

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 82
